$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.678.46'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '3.682.62'
$ws.Range("E3").Value = '  +3.06%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.67'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.88'
$ws.Range("E6").Value = '  +15.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '666.83'
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.428'
$ws.Range("E8").Value = '  +4.84%  '
$ws.Range("E9").Value = '  +4.54%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '3.680.19'
$ws.Range("E11").Value = '  +3.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.48'
$ws.Range("E12").Value = '  +4.50%  '
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("E14").Value = '  +3.34%  '
$ws.Range("D15").Value = '4.368.95'
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000268'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").Value = '96.445.78'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.08'
$ws.Range("E18").Value = '  +17.05%  '
$ws.Range("D19").Value = '3.688.71'
$ws.Range("E19").Value = '  +3.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.89'
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.44'
$ws.Range("E21").Value = '  +3.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.535'
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '518.32'
$ws.Range("E23").Value = '  +3.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.48'
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("E25").Value = '  +4.32%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.62'
$ws.Range("E27").Value = '  +6.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.02'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.167'
$ws.Range("E29").Value = '  +9.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.05'
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.15'
$ws.Range("E31").Value = '  +6.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.90'
$ws.Range("E34").Value = '  +4.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.79'
$ws.Range("E36").Value = '  +9.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.586'
$ws.Range("E37").Value = '  +3.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '616.99'
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.72'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '42.46'
$ws.Range("E40").Value = '  +27.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.160'
$ws.Range("E41").Value = '  +6.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.958'
$ws.Range("E42").Value = '  +6.04%  '
$ws.Range("E43").Value = '  +7.01%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.19'
$ws.Range("E45").Value = '  +8.31%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0457'
$ws.Range("E46").Value = '  +7.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.433'
$ws.Range("E47").Value = '  +25.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.62'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.62'
$ws.Range("E50").Value = '  +4.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.64'
$ws.Range("E51").Value = '  +3.69%  '
